$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Update the shared "Git Commit ID" text for every row that references it (AJ2:AJ80)
$ws.Range("AJ2:AJ80").Value = "IndicatorQuantiles.R, Git Commit ID: d77a77d64f72a744c78cd38270c72c5d9c8cd498"

# Update the pid column (AH2:AH80) from 19980 to 11992
$ws.Range("AH2:AH80").Value = 11992
